$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 11: Asus Radeon R9 290X 4GB Video Card
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "Asus Radeon R9 290X 4GB Video Card"
$ws.Range("C11").Value = "5.0 Average "
$ws.Range("D11").Value = "Asus"
$ws.Range("E11").Value = "R9290X-4GD5"
$ws.Range("F11").Value = "PCI-Express x16"
$ws.Range("G11").Value = "Radeon R9 290X"
$ws.Range("H11").Value = "4GB"
$ws.Range("I11").Value = "GDDR5"
$ws.Range("J11").Value = "1.0GHz"
$ws.Range("K11").Value = "300 Watts"
$ws.Range("L11").Value = "yes"
$ws.Range("M11").Value = "no"
$ws.Range("N11").Value = "4-way CrossFire"
$ws.Range("O11").Value = "no"
$ws.Range("P11").Value = "10.87`" (276mm)"
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = "yes "
$ws.Range("V11").Value = 434.11

# ---------------------------------------------------------------------------
# Row 12: Sapphire Radeon R9 290 4GB Vapor-X Video Card
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "Sapphire Radeon R9 290 4GB Vapor-X Video Card"
$ws.Range("C12").Value = "4.8 Average"
$ws.Range("D12").Value = "Sapphire"
$ws.Range("E12").Value = "100362VXSR"
$ws.Range("F12").Value = "PCI-Express x16"
$ws.Range("G12").Value = "Radeon R9 290"
$ws.Range("H12").Value = "4GB"
$ws.Range("I12").Value = "GDDR5"
$ws.Range("J12").Value = "1.03Ghz"
$ws.Range("K12").Value = "250 Watts"
$ws.Range("L12").Value = "yes"
$ws.Range("M12").Value = "no"
$ws.Range("N12").Value = "4-way CrossFire"
$ws.Range("O12").Value = "no"
$ws.Range("P12").Value = "12.01`" (305mm)"
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 0
$ws.Range("S12").Value = 1
$ws.Range("T12").Value = 1
$ws.Range("U12").Value = "yes "
$ws.Range("V12").Value = 424.3

# ---------------------------------------------------------------------------
# Match the row-15.75pt / wrapped-column look used by rows 5-10 directly
# above: copy the cell formatting down from row 8 (closest template row)
# and patch the couple of columns whose format differs from that template.
# ---------------------------------------------------------------------------
$ws.Range("B8:V8").Copy()
$ws.Range("B11:V11").PasteSpecial(-4122)
$ws.Range("B8:V8").Copy()
$ws.Range("B12:V12").PasteSpecial(-4122)

$ws.Range("H1").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H12").PasteSpecial(-4122)

$ws.Range("C9").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 15.75

$ws.Range("B13").Select()
